# Update InputData/cpi.xlsx to version 3.2.1
#  - regenerate the "Year" row labels in sheet "Data" (A6:A58) using the
#    new label spacing (77 dots + 1 space + period, instead of 77 dots +
#    5 spaces + period), which also re-packs the shared-string table
#  - append the new 2020 data row (row 58) to the "Data" sheet
#  - leave the "About" sheet's URL / year cells alone (their shared
#    string index shifts automatically because of the re-pack above)
#  - make "Data" the active tab, with B58 selected, and leave "About"
#    with B6 selected (its hyperlink cell)

$wb   = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$data  = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# 1. Rewrite the year labels in column A (rows 6-57) with the new
#    spacing, and create row 58 for 2020.
# ---------------------------------------------------------------------
$dots = ''.PadLeft(77, '.')
for ($year = 1968; $year -le 2019; $year++) {
    $r = $year - 1962
    $data.Cells.Item($r, 1).Value = [string]$year + $dots + ' .'
}

# ---------------------------------------------------------------------
# 2. Append the new 2020 row of CPI data.
# ---------------------------------------------------------------------
$data.Cells.Item(58, 1).Value = '2020' + $dots + ' .'
$data.Cells.Item(58, 2).Value = 257.55700000000002
$data.Cells.Item(58, 3).Value = 260.065
$data.Cells.Item(58, 4).Value = 258.81099999999998
$data.Cells.Item(58, 5).Value = 1.4
$data.Cells.Item(58, 6).Value = 1.2
$data.Cells.Item(58, 7).NumberFormat = $data.Cells.Item(57, 7).NumberFormat
$data.Cells.Item(58, 7).Formula = '=$D$50/D58'

# ---------------------------------------------------------------------
# 3. Selection / active sheet bookkeeping, matching the saved view
#    state: "Data" becomes the active tab (scrolled down near the new
#    row) with B58 selected; "About" keeps B6 (its hyperlink) selected.
# ---------------------------------------------------------------------
$about.Activate()
$about.Range("B6").Select()

$data.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$data.Range("B58").Select()
